$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'90.774.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.75%  '

# Row 3
$ws.Range("D3").Value = "'3.159.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.67%  '

# Row 4
$ws.Range("E4").Value = '  +0.24%  '

# Row 5
$ws.Range("D5").Value = "'215.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.43%  '

# Row 6
$ws.Range("D6").Value = "'625.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.99%  '

# Row 7
$ws.Range("D7").Value = "'1.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +26.28%  '

# Row 8
$ws.Range("D8").Value = "'0.369"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.94%  '

# Row 9
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.00%  '

# Row 10
$ws.Range("D10").Value = "'3.157.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.75%  '

# Row 11
$ws.Range("D11").Value = "'0.751"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.97%  '

# Row 12
$ws.Range("E12").Value = '  +6.14%  '

# Row 13
$ws.Range("D13").Value = "'5.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.89%  '

# Row 14
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.20%  '

# Row 15
$ws.Range("D15").Value = "'35.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.20%  '

# Row 16
$ws.Range("D16").Value = "'90.625.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.50%  '

# Row 17
$ws.Range("D17").Value = "'3.741.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.06%  '

# Row 18
$ws.Range("D18").Value = "'3.167.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.23%  '

# Row 19
$ws.Range("D19").Value = "'3.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.42%  '

# Row 20
$ws.Range("D20").Value = "'14.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.07%  '

# Row 21
$ws.Range("D21").Value = "'0.0000214"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.03%  '

# Row 22
$ws.Range("D22").Value = "'469.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.07%  '

# Row 23
$ws.Range("D23").Value = "'9.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.58%  '

# Row 24
$ws.Range("D24").Value = "'5.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.72%  '

# Row 25
$ws.Range("D25").Value = "'96.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.89%  '

# Row 26
$ws.Range("D26").Value = "'5.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.44%  '

# Row 27
$ws.Range("D27").Value = "'12.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.00%  '

# Row 28
$ws.Range("D28").Value = "'3.320.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.16%  '

# Row 29
$ws.Range("E29").Value = '  +0.05%  '

# Row 30
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'0.225"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +58.01%  '

# Row 31
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = "'0.164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.80%  '

# Row 32
$ws.Range("D32").Value = "'9.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.19%  '

# Row 33
$ws.Range("E33").Value = '  -1.37%  '

# Row 34
$ws.Range("D34").Value = "'27.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +17.86%  '

# Row 35
$ws.Range("D35").Value = "'521.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.47%  '

# Row 36
$ws.Range("E36").Value = '  +5.53%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = "'0.145"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.43%  '

# Row 38
$ws.Range("D38").Value = "'7.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.02%  '

# Row 39
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = "'1.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.23%  '

# Row 40
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = "'3.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.28%  '

# Row 41
$ws.Range("D41").Value = "'0.0913"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +26.05%  '

# Row 42
$ws.Range("D42").Value = "'0.431"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +16.04%  '

# Row 43
$ws.Range("D43").Value = "'22.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.32%  '

# Row 44
$ws.Range("E44").Value = '  +0.01%  '

# Row 45
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = "'0.752"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +23.03%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = "'1.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.92%  '

# Row 48
$ws.Range("D48").Value = "'4.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.88%  '

# Row 49
$ws.Range("D49").Value = "'150.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.07%  '

# Row 50
$ws.Range("D50").Value = "'1.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.96%  '

# Row 51
$ws.Range("D51").Value = "'45.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.23%  '
